$wb = $excel.ActiveWorkbook

# Rename the existing sheet from "Sheet 1" to "Invoice"
$invoice = $wb.Worksheets.Item("Sheet 1")
$invoice.Name = "Invoice"

# Add a new worksheet called "Customer" after the Invoice sheet
$customer = $wb.Worksheets.Add($null, $invoice)
$customer.Name = "Customer"

# Populate the Customer sheet with data
$customer.Range("A1").Value = "customer_id"
$customer.Range("B1").Value = "customer_name"
$customer.Range("A2").Value = 2031
$customer.Range("B2").Value = "Ted Zelinsky"

# Select B2 and make Customer the active sheet/tab
$customer.Range("B2").Select()
$customer.Activate()
